$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value2 = '51.825.10'
$ws.Range('E2').Value2 = '  -0.71%  '
$ws.Range('D3').Value2 = '2.778.65'
$ws.Range('E3').Value2 = '  -2.20%  '
$ws.Range('D4').Value = "'1.00"
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value2 = '  -0.08%  '
$ws.Range('D5').Value = "'359.29"
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value2 = '  -0.80%  '
$ws.Range('D6').Value = "'108.82"
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value2 = '  -5.98%  '
$ws.Range('D7').Value = "'0.557"
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value2 = '  +1.02%  '
$ws.Range('E8').Value2 = '  -0.01%  '
$ws.Range('D9').Value = "'0.590"
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value2 = '  -3.18%  '
$ws.Range('E10').Value2 = '  -5.90%  '
$ws.Range('D11').Value = "'0.0848"
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value2 = '  -1.83%  '
$ws.Range('E12').Value2 = '  +0.01%  '
$ws.Range('D13').Value = "'19.41"
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value2 = '  -3.67%  '
$ws.Range('D14').Value = "'7.61"
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value2 = '  -3.41%  '
$ws.Range('D15').Value2 = '3.214.06'
$ws.Range('E15').Value2 = '  -2.47%  '
$ws.Range('D16').Value2 = '2.798.07'
$ws.Range('E16').Value2 = '  -2.34%  '
$ws.Range('D17').Value = "'0.914"
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value2 = '  +1.13%  '
$ws.Range('D18').Value2 = '51.667.55'
$ws.Range('E18').Value2 = '  -1.09%  '
$ws.Range('E19').Value2 = '  +1.01%  '
$ws.Range('D20').Value = "'3.10"
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value2 = '  -2.54%  '
$ws.Range('D21').Value = "'13.04"
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value2 = '  -6.58%  '
$ws.Range('D22').Value2 = '0.0₃0978'
$ws.Range('E22').Value2 = '  -1.73%  '
$ws.Range('D23').Value = "'273.62"
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value2 = '  +0.97%  '
$ws.Range('D24').Value = "'69.38"
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value2 = '  -1.53%  '
$ws.Range('D25').Value = "'2.75"
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value2 = '  -3.09%  '
$ws.Range('D26').Value = "'26.41"
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value2 = '  -2.85%  '
$ws.Range('E27').Value2 = '  +0.11%  '
$ws.Range('D28').Value = "'10.14"
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value2 = '  -1.68%  '
$ws.Range('E29').Value2 = '  -0.99%  '
$ws.Range('E30').Value2 = '  +0.69%  '
$ws.Range('D31').Value = "'0.0469"
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value2 = '  +6.50%  '
$ws.Range('D32').Value = "'51.62"
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value2 = '  +0.89%  '
$ws.Range('D33').Value = "'34.13"
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value2 = '  -1.29%  '
$ws.Range('D34').Value = "'5.72"
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value2 = '  -2.15%  '
$ws.Range('D35').Value = "'5.35"
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value2 = '  +7.15%  '
$ws.Range('D36').Value = "'0.0836"
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value2 = '  -0.02%  '
$ws.Range('E37').Value2 = '  -0.16%  '
$ws.Range('D38').Value = "'3.17"
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value2 = '  -3.58%  '
$ws.Range('E39').Value2 = '  -6.41%  '
$ws.Range('D40').Value = "'17.97"
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value2 = '  -3.73%  '
$ws.Range('E41').Value2 = '  -1.71%  '
$ws.Range('D42').Value = "'125.25"
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value2 = '  -2.09%  '
$ws.Range('D43').Value = "'2.51"
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value2 = '  -3.71%  '
$ws.Range('D44').Value = "'2.24"
$ws.Range('D44').ClearFormats()
$ws.Range('D45').Value = "'21.79"
$ws.Range('D45').ClearFormats()
$ws.Range('D46').Value2 = '2.052.41'
$ws.Range('E46').Value2 = '  -1.10%  '
$ws.Range('E47').Value2 = '  +0.48%  '
$ws.Range('E48').Value2 = '  -5.53%  '
$ws.Range('D49').Value = "'5.70"
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value2 = '  +1.31%  '
$ws.Range('D50').Value = "'0.929"
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value2 = '  -2.28%  '
$ws.Range('D51').Value = "'8.93"
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value2 = '  -0.94%  '
